$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1893.25
$ws.Range("J17").Value = 1893.25
$ws.Range("L17").Value = 5679.75
$ws.Range("N17").Value = -6015.75

$ws.Range("H50").Value = 200
$ws.Range("J50").Value = 200
$ws.Range("L50").Value = 600
$ws.Range("N50").Value = -1550

$ws.Range("H74").Value = 8766.532999999999
$ws.Range("I74").Value = 6214
$ws.Range("K74").Value = 6214
$ws.Range("M74").Value = -5278

$ws.Range("H77").Value = 8766.532999999999
$ws.Range("I77").Value = 6214
$ws.Range("K77").Value = 31070
$ws.Range("M77").Value = -26390

$ws.Range("H112").Value = 1543.4783
$ws.Range("J112").Value = 1570.5349
$ws.Range("L112").Value = 4711.6047
$ws.Range("N112").Value = -6927.6047

$ws.Range("H137").Value = 1734.8966
$ws.Range("I137").Value = 1245.1666
$ws.Range("J137").Value = 2536.2727
$ws.Range("K137").Value = 3735.4998
$ws.Range("L137").Value = 7608.8181
$ws.Range("M137").Value = -1185.4998
$ws.Range("N137").Value = -12708.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13935.458
$ws.Range("I2").Value = 14938.682
$ws.Range("J2").Value = 2900
$ws.Range("K2").Value = 14938.682
$ws.Range("L2").Value = 2900
$ws.Range("M2").Value = -14825.682
$ws.Range("N2").Value = -3126

$ws.Range("H45").Value = 1948.1
$ws.Range("I45").Value = 1290.0769
$ws.Range("K45").Value = 1290.0769
$ws.Range("M45").Value = -913.0769

$ws.Range("H61").Value = 3808.9
$ws.Range("I61").Value = 1739.9
$ws.Range("K61").Value = 1739.9
$ws.Range("M61").Value = -1527.9

$ws.Range("H63").Value = 5389.375
$ws.Range("J63").Value = 2166.3333
$ws.Range("L63").Value = 2166.3333
$ws.Range("N63").Value = -3538.3333

$ws.Range("H66").Value = 5389.375
$ws.Range("J66").Value = 2166.3333
$ws.Range("L66").Value = 10831.6665
$ws.Range("N66").Value = -17695.6665

$ws.Range("H116").Value = 13935.458
$ws.Range("I116").Value = 14938.682
$ws.Range("J116").Value = 2900
$ws.Range("K116").Value = 14938.682
$ws.Range("L116").Value = 2900
$ws.Range("M116").Value = -12644.682
$ws.Range("N116").Value = -7488

$ws.Range("H135").Value = 61399.066
$ws.Range("J135").Value = 61399.066
$ws.Range("L135").Value = 61399.066
$ws.Range("N135").Value = -71539.06599999999

$ws.Range("H136").Value = 3808.9
$ws.Range("I136").Value = 1739.9
$ws.Range("K136").Value = 5219.700000000001
$ws.Range("M136").Value = -2669.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13935.458
$ws.Range("I3").Value = 14938.682
$ws.Range("J3").Value = 2900
$ws.Range("K3").Value = 14938.682
$ws.Range("L3").Value = 2900
$ws.Range("M3").Value = -14824.682
$ws.Range("N3").Value = -3128

$ws.Range("H86").Value = 4328.5835
$ws.Range("I86").Value = 2121.1538
$ws.Range("J86").Value = 6937.364
$ws.Range("K86").Value = 2121.1538
$ws.Range("L86").Value = 6937.364
$ws.Range("M86").Value = -998.1538
$ws.Range("N86").Value = -9183.364

$ws.Range("H89").Value = 4328.5835
$ws.Range("I89").Value = 2121.1538
$ws.Range("J89").Value = 6937.364
$ws.Range("K89").Value = 10605.769
$ws.Range("L89").Value = 34686.82
$ws.Range("M89").Value = -4989.769
$ws.Range("N89").Value = -45918.82

$ws.Range("H94").Value = 3681.0588
$ws.Range("I94").Value = 2971.8667
$ws.Range("K94").Value = 2971.8667
$ws.Range("M94").Value = -2520.8667

$ws.Range("H105").Value = 2679
$ws.Range("I105").Value = 2679
$ws.Range("K105").Value = 2679
$ws.Range("M105").Value = -932

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4111.1055
$ws.Range("I31").Value = 1432
$ws.Range("K31").Value = 1432
$ws.Range("M31").Value = -1137

$ws.Range("H34").Value = 4111.1055
$ws.Range("I34").Value = 1432
$ws.Range("K34").Value = 1432
$ws.Range("M34").Value = -1230

$ws.Range("H58").Value = 2727
$ws.Range("I58").Value = 2332.9092
$ws.Range("K58").Value = 2332.9092
$ws.Range("M58").Value = -2129.9092

$ws.Range("H109").Value = 48904.418
$ws.Range("J109").Value = 48904.418
$ws.Range("L109").Value = 48904.418
$ws.Range("N109").Value = -50984.418

$ws.Range("H132").Value = 1911.2046
$ws.Range("I132").Value = 1831.683
$ws.Range("J132").Value = 2998
$ws.Range("K132").Value = 5495.049
$ws.Range("L132").Value = 8994
$ws.Range("M132").Value = -2965.049
$ws.Range("N132").Value = -14054

$ws.Range("H134").Value = 3202.8147
$ws.Range("I134").Value = 1860.919
$ws.Range("J134").Value = 6123.4116
$ws.Range("K134").Value = 5582.757000000001
$ws.Range("L134").Value = 18370.2348
$ws.Range("M134").Value = -3047.757000000001
$ws.Range("N134").Value = -23440.2348

$ws.Range("H136").Value = 2727
$ws.Range("I136").Value = 2332.9092
$ws.Range("K136").Value = 6998.7276
$ws.Range("M136").Value = -4448.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 87.14286
$ws.Range("I33").Value = 73.333336
$ws.Range("K33").Value = 440.000016
$ws.Range("M33").Value = -157.000016

$ws.Range("H56").Value = 7054.5
$ws.Range("I56").Value = 7054.5
$ws.Range("K56").Value = 7054.5
$ws.Range("M56").Value = -6524.5

$ws.Range("H119").Value = 5932.933
$ws.Range("I119").Value = 1332.3334
$ws.Range("K119").Value = 3997.0002
$ws.Range("M119").Value = 840.9998000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5274
$ws.Range("I126").Value = 3435
$ws.Range("K126").Value = 10305
$ws.Range("M126").Value = -7835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3373.8518
$ws.Range("I7").Value = 2057.7222
$ws.Range("K7").Value = 2057.7222
$ws.Range("M7").Value = -1945.7222

$ws.Range("H40").Value = 11180.667
$ws.Range("J40").Value = 7718.636
$ws.Range("L40").Value = 7718.636
$ws.Range("N40").Value = -7990.636

$ws.Range("H55").Value = 1244
$ws.Range("I55").Value = 1509.4286
$ws.Range("J55").Value = 624.6667
$ws.Range("K55").Value = 1509.4286
$ws.Range("L55").Value = 624.6667
$ws.Range("M55").Value = -1336.4286
$ws.Range("N55").Value = -970.6667

$ws.Range("H68").Value = 5459.048
$ws.Range("I68").Value = 5168.625
$ws.Range("J68").Value = 5637.769
$ws.Range("K68").Value = 5168.625
$ws.Range("L68").Value = 5637.769
$ws.Range("M68").Value = -4419.625
$ws.Range("N68").Value = -7135.769

$ws.Range("H71").Value = 5459.048
$ws.Range("I71").Value = 5168.625
$ws.Range("J71").Value = 5637.769
$ws.Range("K71").Value = 25843.125
$ws.Range("L71").Value = 28188.845
$ws.Range("M71").Value = -22099.125
$ws.Range("N71").Value = -35676.845

$ws.Range("H100").Value = 84542.78999999999
$ws.Range("I100").Value = 280776
$ws.Range("K100").Value = 280776
$ws.Range("M100").Value = -280235

$ws.Range("H126").Value = 3373.8518
$ws.Range("I126").Value = 2057.7222
$ws.Range("K126").Value = 6173.1666
$ws.Range("M126").Value = -3703.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 28571.428
$ws.Range("J109").Value = 28571.428
$ws.Range("L109").Value = 28571.428
$ws.Range("N109").Value = -31345.428
